$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 818
$wsExpo.Range("F4").Value = 281
$wsExpo.Range("F5").Value = 993
$wsExpo.Range("F6").Value = 2349

# Sheet "全部类型" (All types) - same underlying rows, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 818
$wsAll.Range("F4").Value = 281
$wsAll.Range("F7").Value = 993
$wsAll.Range("F8").Value = 2349
